$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update raw data cells (columns C,D,E,F,G,H) for rows 3-8 ---
# Row 3
$ws.Range("C3").Value = 512
$ws.Range("D3").Value = 220
$ws.Range("E3").Value = 0.4296875
$ws.Range("F3").Value = [double]"1.6307830810546799E-3"
$ws.Range("G3").Value = [double]"5.07116317749023E-4"
$ws.Range("H3").Value = [double]"5.1236152648925705E-4"

# Row 4
$ws.Range("C4").Value = 12288
$ws.Range("D4").Value = 5208
$ws.Range("E4").Value = 0.423828125
$ws.Range("F4").Value = [double]"1.16541385650634E-2"
$ws.Range("G4").Value = [double]"1.11699104309082E-2"
$ws.Range("H4").Value = [double]"1.0612964630126899E-2"

# Row 5
$ws.Range("C5").Value = 196608
$ws.Range("D5").Value = 83024
$ws.Range("E5").Value = 0.42228190104166602
$ws.Range("F5").Value = 0.16847181320190399
$ws.Range("G5").Value = 0.16394782066345201
$ws.Range("H5").Value = 0.167225360870361

# Row 6
$ws.Range("C6").Value = 2621440
$ws.Range("D6").Value = 1111264
$ws.Range("E6").Value = 0.42391357421874998
$ws.Range("F6").Value = 2.2640810012817298
$ws.Range("G6").Value = 2.32328748703002
$ws.Range("H6").Value = 2.3609602451324401

# Row 7
$ws.Range("C7").Value = 31457280
$ws.Range("D7").Value = 13440576
$ws.Range("E7").Value = 0.42726440429687501
$ws.Range("F7").Value = 26.398905515670702
$ws.Range("G7").Value = 26.694002389907801
$ws.Range("H7").Value = 26.438570737838699

# Row 8
$ws.Range("C8").Value = 352321536
$ws.Range("D8").Value = 151969152
$ws.Range("E8").Value = 0.431336539132254
$ws.Range("F8").Value = 329.74270033836302
$ws.Range("G8").Value = 332.94395828247002
$ws.Range("H8").Value = 332.61924409866299

# Row 9: clear the data so the row becomes an empty/trailing bucket
$ws.Range("C9").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("F9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()

# --- Number formats follow Excel's automatic "extra precision" bump that
#     happens when these high-precision values are (re)entered. ---
$ws.Range("E3").NumberFormat = "0.0000000"
$ws.Range("F3").NumberFormat = "0.00000000000000000"
$ws.Range("G3").NumberFormat = "0.000000000000000000"
$ws.Range("H3").NumberFormat = "0.000000000000000000"
$ws.Range("I3").NumberFormat = "0.000000000000000000"

$ws.Range("E4").NumberFormat = "0.00000000"
$ws.Range("F4").NumberFormat = "0.0000000000000000"
$ws.Range("G4").NumberFormat = "0.0000000000000000"
$ws.Range("H4").NumberFormat = "0.0000000000000000"
$ws.Range("I4").NumberFormat = "0.0000000000000000"

$ws.Range("E5").NumberFormat = "0.000000000000000"
$ws.Range("F5").NumberFormat = "0.000000000000000"
$ws.Range("G5").NumberFormat = "0.000000000000000"
$ws.Range("H5").NumberFormat = "0.000000000000000"
$ws.Range("I5").NumberFormat = "0.000000000000000"

$ws.Range("E6").NumberFormat = "0.00000000000000"
$ws.Range("F6").NumberFormat = "0.00000000000000"
$ws.Range("G6").NumberFormat = "0.00000000000000"
$ws.Range("H6").NumberFormat = "0.00000000000000"
$ws.Range("I6").NumberFormat = "0.00000000000000"

$ws.Range("E7").NumberFormat = "0.000000000000000"
$ws.Range("F7").NumberFormat = "0.0000000000000"
$ws.Range("G7").NumberFormat = "0.0000000000000"
$ws.Range("H7").NumberFormat = "0.0000000000000"
$ws.Range("I7").NumberFormat = "0.0000000000000"

$ws.Range("E8").NumberFormat = "0.000000000000000"
$ws.Range("F8").NumberFormat = "0.000000000000"
$ws.Range("G8").NumberFormat = "0.00000000000"
$ws.Range("H8").NumberFormat = "0.000000000000"
$ws.Range("I8").NumberFormat = "0.000000000000"

$ws.Range("E9").NumberFormat = "0.000000000000000"
$ws.Range("F9").NumberFormat = "0.00000000000"
$ws.Range("G9").NumberFormat = "0.0000000000"
$ws.Range("H9").NumberFormat = "0.00000000000"
$ws.Range("I9").NumberFormat = "0.00000000000"

# --- Sheet view changes ---
$window = $ws.Application.ActiveWindow
$window.ScrollColumn = 5
$window.ScrollRow = 1
$window.Zoom = 121
$ws.Range("N12").Select()

$excel.CalculateFull()
